$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Names (column A) ---
$ws.Range("A4").Value = "Lythronax"
$ws.Range("A5").Value = "Mapusaurus"
$ws.Range("A6").Value = "Spinosaurus"

# --- Photo hyperlinks (column C) ---
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/lythronax.png?raw=true")
$ws.Range("C4").Value = "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/lythronax.png?raw=true"

$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/mapusaurus.jpg?raw=true", "", "", "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/carcharodontosaurus.jpg?raw=true")
$ws.Range("C5").Value = "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/mapusaurus.jpg?raw=true"

$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/spinosaurus.jpg?raw=true")
$ws.Range("C6").Value = "https://github.com/SergioGerman412/30DayChartChallenge/blob/main/Day19-Dinosaurs/Imágenes/spinosaurus.jpg?raw=true"

# --- Descriptions (column B) ---
$ws.Range("B4").Value = "Apodado ""rey del oeste sangriento"", fue un dinosaurio carnívoro que vivió en América del Norte durante el período Cretácico Superior, caracterizado por su gran tamaño y sus afilados dientes serrados."
$ws.Range("B5").Value = "Terópodo gigante que vivió en lo que ahora es Argentina durante el período Cretácico. Conocido por ser uno de los mayores carnívoros que haya existido, este dinosaurio cazaba en manadas y se caracterizaba por su gran tamaño y ferocidad."
$ws.Range("B6").Value = "Habitó en lo que ahora es el norte de África durante el período Cretácico. Es uno de los mayores dinosaurios carnívoros conocidos"

# --- Re-apply the existing hyperlink cell style so the new cells match C2/C3 ---
$ws.Range("C4").Style = "Hipervínculo"
$ws.Range("C5").Style = "Hipervínculo"
$ws.Range("C6").Style = "Hipervínculo"

# --- Selection as left by the author ---
[void]$ws.Range("B12").Select()
